# Auto-generated edit script applying the Cactuar_Profits.xlsx scheduled-runner update.
# For each affected row, set new currentAveragePrice / LevePrice / LeveProfit values
# exactly as computed by the upstream price-refresh job.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 48
$ws.Range("H48").Value = 6430.231
$ws.Range("J48").Value = 8182.1665
$ws.Range("L48").Value = 24546.4995
$ws.Range("N48").Value = -25130.4995
# Row 56
$ws.Range("H56").Value = 6430.231
$ws.Range("J56").Value = 8182.1665
$ws.Range("L56").Value = 24546.4995
$ws.Range("N56").Value = -25614.4995
# Row 94
$ws.Range("H94").Value = 13674.75
$ws.Range("I94").Value = 7352.5
$ws.Range("K94").Value = 7352.5
$ws.Range("M94").Value = -6901.5
# Row 132
$ws.Range("H132").Value = 4913.2334
$ws.Range("I132").Value = 1676.52
$ws.Range("J132").Value = 21096.8
$ws.Range("K132").Value = 5029.559999999999
$ws.Range("L132").Value = 63290.39999999999
$ws.Range("M132").Value = -2499.559999999999
$ws.Range("N132").Value = -68350.39999999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 16338.594
$ws.Range("I2").Value = 18600.428
$ws.Range("J2").Value = 505.75
$ws.Range("K2").Value = 18600.428
$ws.Range("L2").Value = 505.75
$ws.Range("M2").Value = -18487.428
$ws.Range("N2").Value = -731.75
# Row 45
$ws.Range("H45").Value = 3614.8572
$ws.Range("I45").Value = 2000
$ws.Range("J45").Value = 3884
$ws.Range("K45").Value = 2000
$ws.Range("L45").Value = 3884
$ws.Range("M45").Value = -1623
$ws.Range("N45").Value = -4638
# Row 61
$ws.Range("H61").Value = 4549.727
$ws.Range("I61").Value = 4744.7
$ws.Range("K61").Value = 4744.7
$ws.Range("M61").Value = -4532.7
# Row 98
$ws.Range("H98").Value = 39999
$ws.Range("J98").Value = 39999
$ws.Range("L98").Value = 39999
$ws.Range("N98").Value = -45989
# Row 102
$ws.Range("H102").Value = 2833.8333
$ws.Range("I102").Value = 2833.8333
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2833.8333
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1211.8333
$ws.Range("N102").ClearContents()
# Row 116
$ws.Range("H116").Value = 16338.594
$ws.Range("I116").Value = 18600.428
$ws.Range("J116").Value = 505.75
$ws.Range("K116").Value = 18600.428
$ws.Range("L116").Value = 505.75
$ws.Range("M116").Value = -16306.428
$ws.Range("N116").Value = -5093.75
# Row 122
$ws.Range("H122").Value = 3713.7073
$ws.Range("I122").Value = 3626.1538
$ws.Range("K122").Value = 10878.4614
$ws.Range("M122").Value = -8428.4614
# Row 136
$ws.Range("H136").Value = 4549.727
$ws.Range("I136").Value = 4744.7
$ws.Range("K136").Value = 14234.1
$ws.Range("M136").Value = -11684.1

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 16338.594
$ws.Range("I3").Value = 18600.428
$ws.Range("J3").Value = 505.75
$ws.Range("K3").Value = 18600.428
$ws.Range("L3").Value = 505.75
$ws.Range("M3").Value = -18486.428
$ws.Range("N3").Value = -733.75
# Row 86
$ws.Range("H86").Value = 2592.4583
$ws.Range("I86").Value = 2393
$ws.Range("J86").Value = 2791.9167
$ws.Range("K86").Value = 2393
$ws.Range("L86").Value = 2791.9167
$ws.Range("M86").Value = -1270
$ws.Range("N86").Value = -5037.9167
# Row 89
$ws.Range("H89").Value = 2592.4583
$ws.Range("I89").Value = 2393
$ws.Range("J89").Value = 2791.9167
$ws.Range("K89").Value = 11965
$ws.Range("L89").Value = 13959.5835
$ws.Range("M89").Value = -6349
$ws.Range("N89").Value = -25191.5835
# Row 94
$ws.Range("H94").Value = 1049.1428
$ws.Range("I94").Value = 891
$ws.Range("K94").Value = 891
$ws.Range("M94").Value = -440
# Row 100
$ws.Range("H100").Value = 33299.8
$ws.Range("J100").Value = 33299.8
$ws.Range("L100").Value = 33299.8
$ws.Range("N100").Value = -35463.8
# Row 105
$ws.Range("H105").Value = 4752.0557
$ws.Range("I105").Value = 4601.75
$ws.Range("K105").Value = 4601.75
$ws.Range("M105").Value = -2854.75
# Row 107
$ws.Range("H107").Value = 2175.9285
$ws.Range("I107").Value = 1976.7188
$ws.Range("K107").Value = 1976.7188
$ws.Range("M107").Value = -56.7188000000001
# Row 134
$ws.Range("H134").Value = 2440.1904
$ws.Range("I134").Value = 1703.4814
$ws.Range("K134").Value = 5110.4442
$ws.Range("M134").Value = -2575.4442

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 15387142
$ws.Range("I31").Value = 17859164
$ws.Range("K31").Value = 17859164
$ws.Range("M31").Value = -17858869
# Row 34
$ws.Range("H34").Value = 15387142
$ws.Range("I34").Value = 17859164
$ws.Range("K34").Value = 17859164
$ws.Range("M34").Value = -17858962
# Row 105
$ws.Range("H105").Value = 1262
$ws.Range("I105").Value = 1044.625
$ws.Range("K105").Value = 1044.625
$ws.Range("M105").Value = 702.375
# Row 132
$ws.Range("H132").Value = 53338708
$ws.Range("I132").Value = 60607536
$ws.Range("J132").Value = 33966
$ws.Range("K132").Value = 181822608
$ws.Range("L132").Value = 101898
$ws.Range("M132").Value = -181820078
$ws.Range("N132").Value = -106958

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 52
$ws.Range("H52").Value = 5845.875
$ws.Range("J52").Value = 5845.875
$ws.Range("L52").Value = 17537.625
$ws.Range("N52").Value = -18069.625
# Row 63
$ws.Range("H63").Value = 5817.125
$ws.Range("I63").Value = 2984.25
$ws.Range("K63").Value = 8952.75
$ws.Range("M63").Value = -8203.75
# Row 64
$ws.Range("H64").Value = 8689.125
$ws.Range("I64").Value = 6491.5
$ws.Range("K64").Value = 19474.5
$ws.Range("M64").Value = -19204.5
# Row 66
$ws.Range("H66").Value = 5817.125
$ws.Range("I66").Value = 2984.25
$ws.Range("K66").Value = 26858.25
$ws.Range("M66").Value = -23114.25
# Row 67
$ws.Range("H67").Value = 8689.125
$ws.Range("I67").Value = 6491.5
$ws.Range("K67").Value = 19474.5
$ws.Range("M67").Value = -18538.5
# Row 140
$ws.Range("H140").Value = 3530.3076
$ws.Range("I140").Value = 3361.5833
$ws.Range("K140").Value = 10084.7499
$ws.Range("M140").Value = -4904.749899999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 109472
$ws.Range("I80").Value = 172665.5
$ws.Range("K80").Value = 172665.5
$ws.Range("M80").Value = -171667.5
# Row 83
$ws.Range("H83").Value = 109472
$ws.Range("I83").Value = 172665.5
$ws.Range("K83").Value = 863327.5
$ws.Range("M83").Value = -858335.5
# Row 97
$ws.Range("H97").Value = 1880.6842
$ws.Range("J97").Value = 2492.7778
$ws.Range("L97").Value = 2492.7778
$ws.Range("N97").Value = -3484.7778
# Row 118
$ws.Range("H118").Value = 35000
$ws.Range("J118").Value = 35000
$ws.Range("L118").Value = 35000
$ws.Range("N118").Value = -38314

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 8503.462
$ws.Range("I122").Value = 4509.0713
$ws.Range("J122").Value = 13163.583
$ws.Range("K122").Value = 13527.2139
$ws.Range("L122").Value = 39490.749
$ws.Range("M122").Value = -11077.2139
$ws.Range("N122").Value = -44390.749

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 654.9
$ws.Range("J100").Value = 249.5
$ws.Range("L100").Value = 499
$ws.Range("N100").Value = -1581
# Row 126
$ws.Range("H126").Value = 2616.9092
$ws.Range("I126").Value = 2482.5
$ws.Range("J126").Value = 2778.2
$ws.Range("K126").Value = 7447.5
$ws.Range("L126").Value = 8334.599999999999
$ws.Range("M126").Value = -4977.5
$ws.Range("N126").Value = -13274.6
# Row 138
$ws.Range("H138").Value = 94900
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

